$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update translation strings (D23, D24) to their shortened forms
$ws.Range("D23").Value = "cite the dataset"
$ws.Range("D24").Value = "cite this entry"

# Update the active selection on the sheet (bottom pane) to D26
$ws.Range("D26").Select()
